$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.056.91"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "3.750.09"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.35"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "3.747.40"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.08"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "4.380.55"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "3.744.75"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "69.090.00"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.86"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +17.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  +6.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E31").Value = "  +5.55%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.60"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "3.895.01"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "3.686.12"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.85"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.16"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "2.792.38"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0354"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.35%  "
